$wb = $excel.ActiveWorkbook

# New date for the appended row on every price sheet.
$newDate = "2025-03-31"

# Map of worksheet name -> B-column value to repeat for the new row (same
# value as the prior day, 2025-03-30, per the source diff).
$updates = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "42"
    "N-type Wafer"              = "1.21"
    "Cell Topcon 183mm"         = "0.303"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,533"
    "Silver Busbar front-side"  = "8,284"
    "Silver finger front-side"  = "8,334"
    "USD_CNY"                   = "7.2817"
}

foreach ($name in $updates.Keys) {
    $ws = $wb.Worksheets.Item($name)

    $rowA = $ws.Range("A30")
    $rowB = $ws.Range("B30")

    # Force text storage (matching the existing text-typed date/price cells
    # above them) instead of letting Excel auto-convert to a date serial or
    # a number.
    $rowA.Value = "'" + $newDate
    $rowB.Value = "'" + $updates[$name]

    # Drop the quote-prefix formatting that typing a leading apostrophe
    # applies, so the new cells stay unstyled like their neighbours.
    $ws.Range("A30:B30").ClearFormats()
}
